$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of results for the decision_tree model
$ws.Range("A3").Value = "decision_tree"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.99019607843137203
$ws.Range("D3").Value = 0.95850622406638997
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.66666666666666596
$ws.Range("G3").Value = 0.8
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0.87337662337662303
$ws.Range("J3").Value = 0.62552011095700399

# Match the cursor/selection state recorded after the edit
$ws.Range("J3").Select()
